$d = $word.ActiveDocument

function Find-ParagraphIndexByText($doc, $wantedText, $skipIndex) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($skipIndex -ne $null -and $i -eq $skipIndex) { continue }
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text -eq ($wantedText + "`r")) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. Build the new "Meta description" paragraph next to an already-Normal
#    styled paragraph (so it naturally comes out as a plain, un-styled
#    paragraph with no extra rsid bookkeeping), fill in its text/formatting,
#    then relocate a copy of it to sit right after the document title, and
#    delete the scratch copy we built it in.
# ---------------------------------------------------------------------------
$anchorText = "Are you ready to howl at the moon? Coyote Moon is waiting for you, a 5-reel, 4-row video slot with a whopping 40 adjustable paylines. That’s right, forty! You know what that means? More chances to win big! Woohoo! "
$anchorIdx = Find-ParagraphIndexByText $d $anchorText $null
$anchorPara = $d.Paragraphs.Item($anchorIdx)
$anchorPara.Range.InsertParagraphBefore()

$metaIdx = $anchorIdx
$metaPara = $d.Paragraphs.Item($metaIdx)
$metaText = "Meta description: Discover all you need to know about Coyote Moon slot! Play it for free or for real money at top casinos. Read our full review with the latest information."

$metaRange = $metaPara.Range
$metaRange.MoveEnd(1, -1)
$metaRange.Text = $metaText

$labelLen = ("Meta description").Length
$metaStart = $metaPara.Range.Start
$labelRange = $d.Range($metaStart, $metaStart + $labelLen)
$labelRange.Bold = 1

# Copy (not cut) the finished scratch paragraph, paste it right after the
# title, then delete the scratch paragraph (now shifted down by one since
# the paste added a paragraph above it) - this keeps the relocated
# paragraph free of rsid bookkeeping while avoiding the stray empty
# paragraph that a plain Cut leaves behind in the original spot.
$metaPara.Range.Copy()

$titlePara = $d.Paragraphs.Item(1)
$insertionPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$insertionPoint.Paste()

$scratchPara = $d.Paragraphs.Item($metaIdx + 1)
$scratchPara.Range.Delete()

Write-Host "Step1 done. Paragraph count:" $d.Paragraphs.Count

# ---------------------------------------------------------------------------
# 2. Remove the duplicated bold title paragraph that used to sit near the
#    bottom of the document (right before the italic meta-description
#    paragraph) - its content now lives at the top instead.
# ---------------------------------------------------------------------------
$dupText = "Play Coyote Moon Free: Full Slot Game Review 2021"
$dupIdx = Find-ParagraphIndexByText $d $dupText 1
if ($dupIdx -gt 0) {
    $d.Paragraphs.Item($dupIdx).Range.Delete()
}

Write-Host "Step2 done. Paragraph count:" $d.Paragraphs.Count

# ---------------------------------------------------------------------------
# 3. Replace the text of the remaining italic paragraph (previously the meta
#    description) with the new DALLE image prompt, keeping its italic run
#    formatting and leading empty run intact.
# ---------------------------------------------------------------------------
$oldDescText = "Discover all you need to know about Coyote Moon slot! Play it for free or for real money at top casinos. Read our full review with the latest information."
$dalleText = "Prompt for DALLE: Create a colorful cartoon-style feature image for Coyote Moon that showcases the game's adventurous atmosphere. The image should feature a happy Maya warrior wearing glasses to represent the theme of exploring North America's deserts, forests, and mountains. The warrior should be surrounded by the game's symbols like coyotes, cow skulls, deer, hummingbirds, and lizards and the background should be a golden desert with cacti. The image should be eye-catching and convey the excitement of playing Coyote Moon."

$descIdx = Find-ParagraphIndexByText $d $oldDescText $null
if ($descIdx -gt 0) {
    $r = $d.Paragraphs.Item($descIdx).Range
    $r.MoveEnd(1, -1)
    $r.Text = $dalleText
}

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
